# Apply the corrections to the test model and corresponding test cases.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (CAM_MOUSE_4): steps to execute referenced the wrong test (CAM_MOUSE_2 -> CAM_MOUSE_3)
$ws.Range("D10").Value = "1. Execute test CAM_MOUSE_3`n2. Right click on ARE GUI background panel`n3. Click on 'start button'"

# Row 12 (CAM_MOUSE_6): step 3 should click the 'pause button' instead of 'resume button'
$ws.Range("D12").Value = "1. Execute test CAM_MOUSE_1`n2. Right click on ARE GUI background panel`n3. Click on 'pause button'`n4. Click on 'stop button'"

# Row 13 (CAM_MOUSE_7): requirement link and steps referenced the wrong test (CAM_MOUSE_2 -> CAM_MOUSE_1)
$ws.Range("C13").Value = "Test CAM_MOUSE_1"
$ws.Range("D13").Value = "1. Execute Test CAM_MOUSE_1  by clicking 10 times onto 'Start' button (Play) as fast as possible"

# Row 16 (CAM_MOUSE_10): add the two missing steps to open the webservice URLs
$ws.Range("D16").Value = "1. Open command shell in bin/ARE directory`n2. Execute `nstart.bat --webservice TestModelAutostart.acs`n./start.sh --webservice TestModelAutostart.acs`n3. Open http://localhost:8082/`n4. Open http://localhost:8081/rest/runtime/model`n"

# Adjust row heights to fit the updated content (rows 15 and 16 changed size)
$ws.Rows(15).RowHeight = 81.75
$ws.Rows(16).RowHeight = 192

# Update the saved cursor/selection position
$null = $ws.Range("C16").Select()
